$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.NumberFormat = "General"
    $cell.Style = "Normal"
}

# Row 2
Set-TextValue $ws.Range("D2") "27.124.69"
Set-TextValue $ws.Range("E2") "  -2.60%  "

# Row 3
Set-TextValue $ws.Range("D3") "1.868.77"
Set-TextValue $ws.Range("E3") "  -1.93%  "

# Row 4
Set-TextValue $ws.Range("E4") "  -0.20%  "

# Row 5
Set-TextValue $ws.Range("D5") "307.42"
Set-TextValue $ws.Range("E5") "  -1.89%  "

# Row 6
Set-TextValue $ws.Range("E6") "  -0.24%  "

# Row 7
Set-TextValue $ws.Range("D7") "0.5151"
Set-TextValue $ws.Range("E7") "  +2.32%  "

# Row 8
Set-TextValue $ws.Range("D8") "0.3759"
Set-TextValue $ws.Range("E8") "  -1.56%  "

# Row 9
Set-TextValue $ws.Range("D9") "0.07163"
Set-TextValue $ws.Range("E9") "  -1.56%  "

# Row 10
Set-TextValue $ws.Range("D10") "20.78"
Set-TextValue $ws.Range("E10") "  -0.08%  "

# Row 11
Set-TextValue $ws.Range("D11") "0.8848"
Set-TextValue $ws.Range("E11") "  -2.66%  "

# Row 12
Set-TextValue $ws.Range("D12") "1.879.23"
Set-TextValue $ws.Range("E12") "  -1.58%  "

# Row 13
Set-TextValue $ws.Range("D13") "0.07564"
Set-TextValue $ws.Range("E13") "  -1.16%  "

# Row 14
Set-TextValue $ws.Range("D14") "5.327"
Set-TextValue $ws.Range("E14") "  -2.72%  "

# Row 15
Set-TextValue $ws.Range("D15") "89.21"
Set-TextValue $ws.Range("E15") "  -2.38%  "

# Row 16
Set-TextValue $ws.Range("E16") "  -0.14%  "

# Row 17
Set-TextValue $ws.Range("D17") "0.000008551"
Set-TextValue $ws.Range("E17") "  -1.85%  "

# Row 18
Set-TextValue $ws.Range("D18") "14.18"
Set-TextValue $ws.Range("E18") "  -2.43%  "

# Row 19
Set-TextValue $ws.Range("D19") "1.001"
Set-TextValue $ws.Range("E19") "  -0.25%  "

# Row 20
Set-TextValue $ws.Range("D20") "27.175.76"
Set-TextValue $ws.Range("E20") "  -2.51%  "

# Row 21
Set-TextValue $ws.Range("D21") "5.031"
Set-TextValue $ws.Range("E21") "  -2.62%  "

# Row 22
Set-TextValue $ws.Range("D22") "2.127.33"
Set-TextValue $ws.Range("E22") "  -1.98%  "

# Row 23
Set-TextValue $ws.Range("D23") "10.63"
Set-TextValue $ws.Range("E23") "  -1.70%  "

# Row 24
Set-TextValue $ws.Range("D24") "6.476"
Set-TextValue $ws.Range("E24") "  -1.76%  "

# Row 25
Set-TextValue $ws.Range("D25") "151.51"
Set-TextValue $ws.Range("E25") "  -1.76%  "

# Row 26
Set-TextValue $ws.Range("D26") "1.847"
Set-TextValue $ws.Range("E26") "  -1.13%  "

# Row 27
Set-TextValue $ws.Range("D27") "18.04"
Set-TextValue $ws.Range("E27") "  -1.85%  "

# Row 28
Set-TextValue $ws.Range("D28") "2.166"
Set-TextValue $ws.Range("E28") "  -2.89%  "

# Row 29
Set-TextValue $ws.Range("D29") "112.93"
Set-TextValue $ws.Range("E29") "  -2.01%  "

# Row 30
Set-TextValue $ws.Range("D30") "4.742"
Set-TextValue $ws.Range("E30") "  -3.40%  "

# Row 31
Set-TextValue $ws.Range("D31") "4.691"
Set-TextValue $ws.Range("E31") "  +1.17%  "

# Row 32
Set-TextValue $ws.Range("D32") "0.09022"
Set-TextValue $ws.Range("E32") "  +0.40%  "

# Row 33
Set-TextValue $ws.Range("D33") "0.05159"
Set-TextValue $ws.Range("E33") "  -1.83%  "

# Row 34
Set-TextValue $ws.Range("D34") "3.095"
Set-TextValue $ws.Range("E34") "  -3.37%  "

# Row 35
Set-TextValue $ws.Range("D35") "0.7565"
Set-TextValue $ws.Range("E35") "  -0.99%  "

# Row 36
Set-TextValue $ws.Range("D36") "1.174"
Set-TextValue $ws.Range("E36") "  -4.75%  "

# Row 37
Set-TextValue $ws.Range("D37") "0.02041"
Set-TextValue $ws.Range("E37") "  -0.74%  "

# Row 38
Set-TextValue $ws.Range("D38") "2.512"
Set-TextValue $ws.Range("E38") "  -1.24%  "

# Row 39
Set-TextValue $ws.Range("D39") "3.034"
Set-TextValue $ws.Range("E39") "  +0.54%  "

# Row 40
Set-TextValue $ws.Range("D40") "1.083"
Set-TextValue $ws.Range("E40") "  -1.22%  "

# Row 41
Set-TextValue $ws.Range("D41") "0.5407"
Set-TextValue $ws.Range("E41") "  -2.87%  "

# Row 42
Set-TextValue $ws.Range("D42") "6.670"
Set-TextValue $ws.Range("E42") "  -4.33%  "

# Row 43
Set-TextValue $ws.Range("D43") "114.73"
Set-TextValue $ws.Range("E43") "  +2.85%  "

# Row 44
Set-TextValue $ws.Range("D44") "8.530"
Set-TextValue $ws.Range("E44") "  +0.59%  "

# Row 45
Set-TextValue $ws.Range("D45") "0.1489"
Set-TextValue $ws.Range("E45") "  -1.68%  "

# Row 46
Set-TextValue $ws.Range("D46") "0.4686"
Set-TextValue $ws.Range("E46") "  -2.42%  "

# Row 47
Set-TextValue $ws.Range("B47") "PaxDollar"
Set-TextValue $ws.Range("C47") "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
Set-TextValue $ws.Range("D47") "1.001"
Set-TextValue $ws.Range("E47") "  -0.25%  "

# Row 48
Set-TextValue $ws.Range("B48") "EnergySwap"
Set-TextValue $ws.Range("C48") "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue $ws.Range("D48") "10.16"
Set-TextValue $ws.Range("E48") "  -3.97%  "

# Row 49
Set-TextValue $ws.Range("D49") "1.576"
Set-TextValue $ws.Range("E49") "  -3.35%  "

# Row 50
Set-TextValue $ws.Range("E50") "  -3.46%  "

# Row 51
Set-TextValue $ws.Range("D51") "36.52"
Set-TextValue $ws.Range("E51") "  -1.39%  "
